$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "69.408.39"
$ws.Cells.Item(2, 5).Value = "  -4.06%  "
$ws.Cells.Item(3, 4).Value = "2.505.67"
$ws.Cells.Item(3, 5).Value = "  -5.68%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
$ws.Cells.Item(5, 4).Value = "'579.62"
$ws.Cells.Item(5, 5).Value = "  -2.19%  "
$ws.Cells.Item(6, 4).Value = "'167.39"
$ws.Cells.Item(6, 5).Value = "  -4.17%  "
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
$ws.Cells.Item(8, 4).Value = "'0.520"
$ws.Cells.Item(8, 5).Value = "  -0.38%  "
$ws.Cells.Item(9, 4).Value = "2.505.52"
$ws.Cells.Item(9, 5).Value = "  -5.64%  "
$ws.Cells.Item(10, 4).Value = "'0.158"
$ws.Cells.Item(10, 5).Value = "  -7.54%  "
$ws.Cells.Item(11, 5).Value = "  -0.58%  "
$ws.Cells.Item(12, 5).Value = "  -5.00%  "
$ws.Cells.Item(13, 4).Value = "'4.86"
$ws.Cells.Item(13, 5).Value = "  -2.38%  "
$ws.Cells.Item(14, 4).Value = "2.962.32"
$ws.Cells.Item(14, 5).Value = "  -5.76%  "
$ws.Cells.Item(15, 4).Value = "69.378.66"
$ws.Cells.Item(15, 5).Value = "  -3.89%  "
$ws.Cells.Item(16, 5).Value = "  -5.89%  "
$ws.Cells.Item(17, 5).Value = "  -4.62%  "
$ws.Cells.Item(18, 4).Value = "2.492.29"
$ws.Cells.Item(18, 5).Value = "  -6.18%  "
$ws.Cells.Item(19, 4).Value = "'11.42"
$ws.Cells.Item(19, 5).Value = "  -7.46%  "
$ws.Cells.Item(20, 5).Value = "  -3.39%  "
$ws.Cells.Item(21, 4).Value = "'351.43"
$ws.Cells.Item(21, 5).Value = "  -5.88%  "
$ws.Cells.Item(22, 4).Value = "'3.95"
$ws.Cells.Item(22, 5).Value = "  -5.45%  "
$ws.Cells.Item(23, 5).Value = "  -4.15%  "
$ws.Cells.Item(24, 5).Value = "  -0.08%  "
$ws.Cells.Item(25, 4).Value = "'69.13"
$ws.Cells.Item(25, 5).Value = "  -3.59%  "
$ws.Cells.Item(26, 5).Value = "  -5.65%  "
$ws.Cells.Item(27, 4).Value = "'9.04"
$ws.Cells.Item(27, 5).Value = "  -6.76%  "
$ws.Cells.Item(28, 4).Value = "2.640.62"
$ws.Cells.Item(28, 5).Value = "  -5.61%  "
$ws.Cells.Item(29, 4).Value = "'0.997"
$ws.Cells.Item(29, 5).Value = "  -0.15%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0903"
$ws.Cells.Item(30, 5).Value = "  -5.70%  "
$ws.Cells.Item(31, 5).Value = "  -2.61%  "
$ws.Cells.Item(32, 4).Value = "'480.95"
$ws.Cells.Item(32, 5).Value = "  -3.65%  "
$ws.Cells.Item(33, 5).Value = "  +0.96%  "
$ws.Cells.Item(34, 5).Value = "  -3.05%  "
$ws.Cells.Item(35, 5).Value = "  -0.06%  "
$ws.Cells.Item(36, 5).Value = "  -1.76%  "
$ws.Cells.Item(37, 4).Value = "'152.91"
$ws.Cells.Item(37, 5).Value = "  -5.57%  "
$ws.Cells.Item(38, 4).Value = "'18.87"
$ws.Cells.Item(38, 5).Value = "  -0.15%  "
$ws.Cells.Item(39, 4).Value = "'18.57"
$ws.Cells.Item(39, 5).Value = "  -4.36%  "
$ws.Cells.Item(40, 5).Value = "  +0.00%  "
$ws.Cells.Item(41, 4).Value = "'4.78"
$ws.Cells.Item(41, 5).Value = "  -2.59%  "
$ws.Cells.Item(42, 5).Value = "  -3.04%  "
$ws.Cells.Item(43, 4).Value = "'1.63"
$ws.Cells.Item(43, 5).Value = "  -5.82%  "
$ws.Cells.Item(44, 5).Value = "  -13.57%  "
$ws.Cells.Item(45, 5).Value = "  -8.59%  "
$ws.Cells.Item(46, 4).Value = "'38.18"
$ws.Cells.Item(46, 5).Value = "  -2.52%  "
$ws.Cells.Item(47, 4).Value = "'143.81"
$ws.Cells.Item(47, 5).Value = "  -7.07%  "
$ws.Cells.Item(48, 5).Value = "  -4.01%  "
$ws.Cells.Item(49, 4).Value = "'0.531"
$ws.Cells.Item(49, 5).Value = "  -3.76%  "
$ws.Cells.Item(50, 5).Value = "  -5.15%  "
$ws.Cells.Item(51, 5).Value = "  -2.53%  "
